$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.038.99"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.902.33"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06906"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08028"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7563"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "1.903.98"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.236"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "30.042.91"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007773"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.157.14"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.061"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1266"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.057"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.352"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.041"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7403"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.256"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4466"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.953"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8315"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.688"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.752"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "2.057.39"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1163"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.53%  "
